# Refactor phase1-A notebook: update data loading, enhance markdown sections,
# and correct column names for clarity.
#
# - Add a new "pivot" summary sheet (Offline-votes / E-votes / Total) ahead
#   of the original "Sheet 1" data sheet.
# - Leave "Sheet 1" content intact but move the selection/cursor.
# - Make "Sheet 1" the active tab again (as it was before).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Original data sheet ("Sheet 1") -- just nudge the remembered
#    selection, everything else about it stays the same.
# ---------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Sheet 1")
[void]$dataSheet.Range("A24").Select()

# ---------------------------------------------------------------------
# 2) New "pivot" sheet with the Offline-votes / E-votes / Total summary.
#    Worksheets.Add() inserts the new sheet ahead of the active sheet,
#    i.e. in front of "Sheet 1" -- exactly where it needs to land.
# ---------------------------------------------------------------------
$pivot = $wb.Worksheets.Add()
$pivot.Name = "pivot"

# Re-fetch the data-sheet reference: it was grabbed before the new sheet
# was inserted, and needs to be current for Activate() further down.
$dataSheet = $wb.Worksheets.Item("Sheet 1")

# Header row (same column headers as the data sheet).
$pivot.Range("A1").Value = "party"
$pivot.Range("B1").Value = "Red"
$pivot.Range("C1").Value = "Green"
$pivot.Range("D1").Value = "Invalid ballots"
$pivot.Range("E1").Value = "Total"

# Row 2: Offline-votes (sum via formula).
$pivot.Range("A2").Value = "Offline-votes"
$pivot.Range("B2").Value = 278
$pivot.Range("C2").Value = 406
$pivot.Range("D2").Value = 18
$pivot.Range("E2").Formula = "=SUM(B2:D2)"

# Row 3: E-votes.
$pivot.Range("A3").Value = "E-votes"
$pivot.Range("B3").Value = 130
$pivot.Range("C3").Value = 206
$pivot.Range("D3").Value = 1
$pivot.Range("E3").Value = 337

# Row 4: Total.
$pivot.Range("A4").Value = "Total"
$pivot.Range("B4").Value = 408
$pivot.Range("C4").Value = 612
$pivot.Range("D4").Value = 19
$pivot.Range("E4").Value = 1039

# Column A is wide enough to show the labels, like the data sheet.
$pivot.Columns.Item(1).ColumnWidth = 35.1640625

# Zoom in a bit and remember a selection away from the data, like the author did.
[void]$pivot.Range("I16").Select()
$excel.ActiveWindow.Zoom = 140

# ---------------------------------------------------------------------
# 3) Restore "Sheet 1" as the active/visible tab.
# ---------------------------------------------------------------------
$dataSheet.Activate()
